# Applies the three run-level edits described by the commit diff:
#
#  1. "Index Page Screen Short" paragraph gains a trailing space
#     (as its own run with the paragraph's Times New Roman / sz 36 formatting).
#  2. The "In This Page I can Simply Create HTML file Which is Contents ..."
#     paragraph had its text split across three runs around the word
#     "Which" (wrapped in gramStart/gramEnd proofErr markers) - these are
#     merged back into a single run with no proofErr markers.
#  3. The "In This I can complete Advance Java Assignment Part 1 and in
#     This page I can provide all tasks links  " paragraph - same kind of
#     three-way split around the word "This" - merged back into one run.

$d = $word.ActiveDocument

function Merge-RunsByText($fullText) {
    # Find the full (logical) text even though it currently spans multiple
    # runs / is interrupted by <w:proofErr/> markers, then rewrite it as a
    # single contiguous run. A direct "set Text to the same Text" is a
    # no-op, so we bounce through a short placeholder first to force the
    # engine to actually rebuild the run (merging formatting-identical
    # runs and dropping the now-pointless proofErr markers).
    $rng = $d.Content
    $ok = $rng.Find.Execute($fullText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        return $false
    }
    $startPos = $rng.Start
    $placeholder = "___MERGE_PLACEHOLDER___"
    $rng.Text = $placeholder
    $rng2 = $d.Range($startPos, $startPos + $placeholder.Length)
    $rng2.Text = $fullText
    return $true
}

# --- Change 2: merge the "Which" split across the Index page description ---
$target2 = "In This Page I can Simply Create HTML file Which is Contents Assignments Parts link We can simply Click the Button and we can Redirect Assignment Task Solution Page   "
Merge-RunsByText $target2 | Out-Null

# --- Change 3: merge the "This" split across the Part 1 page description ---
$target3 = "In This I can complete Advance Java Assignment Part 1 and in This page I can provide all tasks links  "
Merge-RunsByText $target3 | Out-Null

# --- Change 1: append a trailing-space run after "Index Page Screen Short" ---
$r1 = $d.Content
$null = $r1.Find.Execute("Index Page Screen Short", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.InsertAfter(" ")
